{"js": "// \"removed Rob from items\"\n// The byline paragraph (\"by Rob Byrd / Chief Enterprise Architect / City of\n// Austin\") is emptied out (only the trailing manual line break survives),\n// and the \"_GoBack\" bookmark that used to sit on the title paragraph is\n// relocated to the start of that now-empty byline paragraph.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph 0 = title (\"Understanding Roles between ... Office\")\n// Paragraph 1 = byline (\"by Rob Byrd / Chief Enterprise Architect / City of Austin\")\nconst bylinePara = paragraphs.items[1];\n\n// Locate the text to remove: everything in the byline paragraph up to (but\n// excluding) the final, lone manual line break that remains afterwards.\nconst byLineBreak = \"\\u000b\"; // manual line break (<w:br/>) as seen in Range/Paragraph text\nconst searchText =\n  \"by Rob Byrd\" + byLineBreak + \"Chief Enterprise Architect\" + byLineBreak + \"City of Austin\";\n\nconst results = bylinePara.search(searchText, { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not locate the 'by Rob Byrd...City of Austin' byline text\");\n}\n\nconst match = results.items[0];\nconst paraStart = bylinePara.getRange(\"Start\");\nconst matchEnd = match.getRange(\"End\");\nconst removalRange = paraStart.expandTo(matchEnd);\n\n// Remove the old \"_GoBack\" bookmark (it currently lives on the title\n// paragraph) before re-inserting it at the byline paragraph's new start.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// Insert the bookmark at the (still valid) start of the byline paragraph,\n// then delete the byline text itself, leaving only the trailing break.\nparaStart.insertBookmark(\"_GoBack\");\nremovalRange.delete();\nawait context.sync();\n", "ps1": "# \"removed Rob from items\"\n# The byline paragraph (\"by Rob Byrd / Chief Enterprise Architect / City of\n# Austin\") is emptied out (only the trailing manual line break survives),\n# and the \"_GoBack\" bookmark that used to sit on the title paragraph is\n# relocated to the start of that now-empty byline paragraph.\n\n$d = $word.ActiveDocument\n\n# Drop the existing \"_GoBack\" bookmark (currently on the title paragraph).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Locate \"by Rob Byrd<br>Chief Enterprise Architect<br>City of Austin\" -\n# i.e. everything in the byline paragraph except the final, lone manual\n# line break that must remain afterwards. Manual line breaks (<w:br/>)\n# show up as Chr(11) in Range.Text.\n$lineBreak = [char]11\n$searchText = \"by Rob Byrd\" + $lineBreak + \"Chief Enterprise Architect\" + $lineBreak + \"City of Austin\"\n\n$found = $d.Content\n$found.Find.ClearFormatting()\n$ok = $found.Find.Execute($searchText)\nif (-not $ok) {\n    throw \"Could not locate the 'by Rob Byrd...City of Austin' byline text\"\n}\n\n$startPos = $found.Start\n\n# Remove the byline text, then re-create \"_GoBack\" as a zero-length\n# bookmark at the same spot (now the start of the byline paragraph,\n# immediately before the remaining line break).\n$found.Delete()\n$bmRange = $d.Range($startPos, $startPos)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n"}
